$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7
$ws.Range("C4").Value = 0.9
$ws.Range("B5").Value = "'FALSE"
$ws.Range("C5").Value = 0.8
$ws.Range("C6").Value = 0.8
$ws.Range("C7").Value = 1
$ws.Range("C8").Value = 0.8
$ws.Range("B9").Value = "'TRUE"
$ws.Range("C9").Value = 1
